$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 24.40667417833229
$ws.Range("C2").Value = 11.68122449503126
$ws.Range("D2").Value = 4.544271183612937
$ws.Range("F2").Value = 52.6222813664527
$ws.Range("G2").Value = 3.763229909519437
$ws.Range("I2").Value = 36.77641339118414
$ws.Range("J2").Value = 9.83704870042666
$ws.Range("L2").Value = 12.76954042613325
$ws.Range("B3").Value = 24.12578674411975
$ws.Range("C3").Value = 11.28937801641274
$ws.Range("D3").Value = 4.519794830385382
$ws.Range("F3").Value = 52.4805796729271
$ws.Range("G3").Value = 3.767544530718476
$ws.Range("I3").Value = 36.75854143372534
$ws.Range("J3").Value = 9.853924868201119
$ws.Range("L3").Value = 12.77687488681967
$ws.Range("B4").Value = 23.95886683829773
$ws.Range("C4").Value = 11.04583855981338
$ws.Range("D4").Value = 4.504441606795146
$ws.Range("F4").Value = 52.4059156209149
$ws.Range("G4").Value = 3.77032972869128
$ws.Range("I4").Value = 36.7551339148644
$ws.Range("J4").Value = 9.864872266130488
$ws.Range("L4").Value = 12.78358451757089
$ws.Range("B5").Value = 23.89232093589902
$ws.Range("C5").Value = 10.9460327889026
$ws.Range("D5").Value = 4.498103236242436
$ws.Range("F5").Value = 52.37859958367894
$ws.Range("G5").Value = 3.771499055580718
$ws.Range("I5").Value = 36.75564238745474
$ws.Range("J5").Value = 9.869481008779692
$ws.Range("L5").Value = 12.78687379674929
$ws.Range("B6").Value = 23.88136244178896
$ws.Range("C6").Value = 10.92943168262093
$ws.Range("D6").Value = 4.49704580824147
$ws.Range("F6").Value = 52.3742517563167
$ws.Range("G6").Value = 3.771695299255912
$ws.Range("I6").Value = 36.75584121577678
$ws.Range("J6").Value = 9.870255212619632
$ws.Range("L6").Value = 12.78745350757941
$ws.Range("B7").Value = 23.95796329930474
$ws.Range("C7").Value = 11.04449457181118
$ws.Range("D7").Value = 4.50435645675542
$ws.Range("F7").Value = 52.40553462812434
$ws.Range("G7").Value = 3.770345359438073
$ws.Range("I7").Value = 36.7551330994167
$ws.Range("J7").Value = 9.864933823141067
$ws.Range("L7").Value = 12.7836266303346
$ws.Range("B8").Value = 24.30872241881221
$ws.Range("C8").Value = 11.54683016067551
$ws.Range("D8").Value = 4.535899014127064
$ws.Range("F8").Value = 52.57086626335769
$ws.Range("G8").Value = 3.764689446132652
$ws.Range("I8").Value = 36.76867786829824
$ws.Range("J8").Value = 9.842746349238411
$ws.Range("L8").Value = 12.77161164914254
$ws.Range("B9").Value = 25.03688615308408
$ws.Range("C9").Value = 12.5013667568965
$ws.Range("D9").Value = 4.595186760771779
$ws.Range("F9").Value = 52.99256067539824
$ws.Range("G9").Value = 3.754671089176333
$ws.Range("I9").Value = 36.85547188427621
$ws.Range("J9").Value = 9.80386334914688
$ws.Range("L9").Value = 12.76553917385487
$ws.Range("B10").Value = 25.59138745281327
$ws.Range("C10").Value = 13.17540945923834
$ws.Range("D10").Value = 4.637194238393945
$ws.Range("F10").Value = 53.3610066011071
$ws.Range("G10").Value = 3.747955900464278
$ws.Range("I10").Value = 36.95617002400056
$ws.Range("J10").Value = 9.778091263532056
$ws.Range("L10").Value = 12.77171000226504
$ws.Range("B11").Value = 25.84678697932002
$ws.Range("C11").Value = 13.47458448345506
$ws.Range("D11").Value = 4.655970435504172
$ws.Range("F11").Value = 53.54111364980231
$ws.Range("G11").Value = 3.745039234142634
$ws.Range("I11").Value = 37.01002034968918
$ws.Range("J11").Value = 9.766968474087824
$ws.Range("L11").Value = 12.77681635291037
$ws.Range("B12").Value = 25.94386209495165
$ws.Range("C12").Value = 13.58669227865317
$ws.Range("D12").Value = 4.663032698571697
$ws.Range("F12").Value = 53.6110872742691
$ws.Range("G12").Value = 3.743954483501737
$ws.Range("I12").Value = 37.03156797287336
$ws.Range("J12").Value = 9.762842589885491
$ws.Range("L12").Value = 12.77907946077821
$ws.Range("B13").Value = 25.92294064746187
$ws.Range("C13").Value = 13.56260213424035
$ws.Range("D13").Value = 4.661513842631376
$ws.Range("F13").Value = 53.59593888930498
$ws.Range("G13").Value = 3.744187228668656
$ws.Range("I13").Value = 37.02687593350785
$ws.Range("J13").Value = 9.763727350328514
$ws.Range("L13").Value = 12.77857742861625
$ws.Range("B14").Value = 25.85476670269317
$ws.Range("C14").Value = 13.48383197029904
$ws.Range("D14").Value = 4.656552411441383
$ws.Range("F14").Value = 53.54683509063852
$ws.Range("G14").Value = 3.744949596438955
$ws.Range("I14").Value = 37.01176993314244
$ws.Range("J14").Value = 9.766627312080683
$ws.Range("L14").Value = 12.77699594501453
$ws.Range("B15").Value = 25.81305245331678
$ws.Range("C15").Value = 13.43542576546379
$ws.Range("D15").Value = 4.653507157787322
$ws.Range("F15").Value = 53.51698739651673
$ws.Range("G15").Value = 3.745419133914865
$ws.Range("I15").Value = 37.00266752990395
$ws.Range("J15").Value = 9.768414821311715
$ws.Range("L15").Value = 12.77607010583883
$ws.Range("B16").Value = 25.57475201667838
$ws.Range("C16").Value = 13.15569788606569
$ws.Range("D16").Value = 4.635960468405218
$ws.Range("F16").Value = 53.34948551362256
$ws.Range("G16").Value = 3.748149279238106
$ws.Range("I16").Value = 36.95281256969787
$ws.Range("J16").Value = 9.778830228280484
$ws.Range("L16").Value = 12.77142245948717
$ws.Range("B17").Value = 25.42930403555936
$ws.Range("C17").Value = 12.98210193857446
$ws.Range("D17").Value = 4.625110831786448
$ws.Range("F17").Value = 53.24991209612024
$ws.Range("G17").Value = 3.749859413385092
$ws.Range("I17").Value = 36.92428780767877
$ws.Range("J17").Value = 9.785373434352691
$ws.Range("L17").Value = 12.76915937961817
$ws.Range("B18").Value = 25.34594859327127
$ws.Range("C18").Value = 12.88155872836335
$ws.Range("D18").Value = 4.618839042471046
$ws.Range("F18").Value = 53.19381835973789
$ws.Range("G18").Value = 3.750856044444601
$ws.Range("I18").Value = 36.90863821898912
$ws.Range("J18").Value = 9.789193509374321
$ws.Range("L18").Value = 12.76807421055815
$ws.Range("B19").Value = 25.31778064727898
$ws.Range("C19").Value = 12.84740070834367
$ws.Range("D19").Value = 4.616710133907564
$ws.Range("F19").Value = 53.17502910385671
$ws.Range("G19").Value = 3.751195724608472
$ws.Range("I19").Value = 36.90346957342193
$ws.Range("J19").Value = 9.790496652469736
$ws.Range("L19").Value = 12.76774400709004
$ws.Range("B20").Value = 25.44475659978808
$ws.Range("C20").Value = 13.00065434750563
$ws.Range("D20").Value = 4.626269039975131
$ws.Range("F20").Value = 53.26039009217816
$ws.Range("G20").Value = 3.749676021455537
$ws.Range("I20").Value = 36.92724596212828
$ws.Range("J20").Value = 9.784671044028299
$ws.Range("L20").Value = 12.76937788994089
$ws.Range("B21").Value = 25.87478198290963
$ws.Range("C21").Value = 13.50700162035237
$ws.Range("D21").Value = 4.658011003114395
$ws.Range("F21").Value = 53.56121022772619
$ws.Range("G21").Value = 3.744725136193211
$ws.Range("I21").Value = 37.01617557940875
$ws.Range("J21").Value = 9.765773189805946
$ws.Range("L21").Value = 12.77745153437449
$ws.Range("B22").Value = 26.15788963653783
$ws.Range("C22").Value = 13.83098855345563
$ws.Range("D22").Value = 4.678477167065448
$ws.Range("F22").Value = 53.76812191154175
$ws.Range("G22").Value = 3.741604372980905
$ws.Range("I22").Value = 37.08103113465827
$ws.Range("J22").Value = 9.753923895480243
$ws.Range("L22").Value = 12.78464770299492
$ws.Range("B23").Value = 26.00663148429303
$ws.Range("C23").Value = 13.65873940671262
$ws.Range("D23").Value = 4.667579512681034
$ws.Range("F23").Value = 53.65675553773556
$ws.Range("G23").Value = 3.743259510775556
$ws.Range("I23").Value = 37.04580084045281
$ws.Range("J23").Value = 9.760202312649142
$ws.Range("L23").Value = 12.78063176750299
$ws.Range("B24").Value = 25.43776966447698
$ws.Range("C24").Value = 12.99226910156335
$ws.Range("D24").Value = 4.625745520546403
$ws.Range("F24").Value = 53.25564939889175
$ws.Range("G24").Value = 3.749758891036388
$ws.Range("I24").Value = 36.92590624537971
$ws.Range("J24").Value = 9.784988413079487
$ws.Range("L24").Value = 12.76927842883855
$ws.Range("B25").Value = 24.83609932650689
$ws.Range("C25").Value = 12.24735515873675
$ws.Range("D25").Value = 4.579421187048572
$ws.Range("F25").Value = 52.86811445682047
$ws.Range("G25").Value = 3.757267368126691
$ws.Range("I25").Value = 36.82551695615766
$ws.Range("J25").Value = 9.813889542774254
$ws.Range("L25").Value = 12.76531198766297
